# Slide 12 ("Chapter 7"), shape "TextBox 13" holds a bulleted paragraph
# describing the decltype fallback rule.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)
$shp = $s.Shapes.Item(9)
$tf = $shp.TextFrame
$tr = $tf.TextRange

# The paragraph currently reads (runs):
#   [1] " "
#   [2] "위에 해당하지 않은 표현식이 "   <- split this run into three
#   [3] "decltype"
#   [4] "에 사용되었을 경우"
#   [5] ", "
#   [6] "아래 방법을 따른다"
#   [7] "."
#
# Split run [2] into three runs, fixing "않은" -> "않는" along the way:
#   "위에 " + "해당하지 않는 " + "표현식이 "
# Using Characters(start, length) (1-based, over the whole TextRange) keeps
# the surrounding runs/formatting untouched while re-writing just that span.
$tr.Characters(2, 3).Text = "위에 "
$tr.Characters(5, 8).Text = "해당하지 않는 "
$tr.Characters(13, 5).Text = "표현식이 "
